$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "L1cam", "Cntn1", "FAPs", 2, 1, 4.495339, 8.990677999999999, 0.1274106381683396, 0.09709078547908702, 3, 1, 0.1145166666666667, 0.34355, 0.5923587361059505, 0.6855055141740312, 0.5147912378166666, 3.0887474269, 0.0754728045918502, 0.0665562688214021),
    @("ECs", "L1cam", "Cntn1", "sCs", 2, 1, 4.495339, 8.990677999999999, 0.1274106381683396, 0.09709078547908702, 1, 0.5, 0.0788065, 0.157613, 0.4076412638940496, 0.3144944858259688, 0.3542619329034999, 1.417047731614, 0.05193783357648938, 0.03053451665768491),
    @("FAPs", "L1cam", "Cntn1", "FAPs", 3, 1, 0.4630963333333333, 1.389289, 0.01312546158663854, 0.01500300202804008, 3, 1, 0.1145166666666667, 0.34355, 0.5923587361059505, 0.6855055141740312, 0.05303224843888889, 0.47729023595, 0.007774981836268409, 0.01028464061938565),
    @("FAPs", "L1cam", "Cntn1", "sCs", 3, 1, 0.4630963333333333, 1.389289, 0.01312546158663854, 0.01500300202804008, 1, 0.5, 0.0788065, 0.157613, 0.4076412638940496, 0.3144944858259688, 0.03649500119283333, 0.218970007157, 0.005350479750370132, 0.004718361408654432),
    @("M1", "L1cam", "Cntn1", "FAPs", 3, 1, 6.798299666666666, 20.394899, 0.1926830655017586, 0.2202455436260365, 3, 1, 0.1145166666666667, 0.34355, 0.5923587361059505, 0.6855055141740312, 0.7785186168277777, 7.00666755145, 0.1141374971496418, 0.1509795346279051),
    @("M1", "L1cam", "Cntn1", "sCs", 3, 1, 6.798299666666666, 20.394899, 0.1926830655017586, 0.2202455436260365, 1, 0.5, 0.0788065, 0.157613, 0.4076412638940496, 0.3144944858259688, 0.5357502026811666, 3.214501216087, 0.07854556835211685, 0.06926600899813132),
    @("M2", "L1cam", "Cntn1", "FAPs", 3, 1, 8.495531999999999, 25.486596, 0.2407874364312792, 0.275231036505607, 3, 1, 0.1145166666666667, 0.34355, 0.5923587361059505, 0.6855055141740312, 0.9728800061999999, 8.755920055800001, 0.1426325415146244, 0.1886723931964276),
    @("M2", "L1cam", "Cntn1", "sCs", 3, 1, 8.495531999999999, 25.486596, 0.2407874364312792, 0.275231036505607, 1, 0.5, 0.0788065, 0.157613, 0.4076412638940496, 0.3144944858259688, 0.6695031425579999, 4.017018855348, 0.09815489491665477, 0.08655864330917931),
    @("Neutro", "L1cam", "Cntn1", "FAPs", 3, 1, 6.279227333333334, 18.837682, 0.1779710855497397, 0.2034290786507157, 3, 1, 0.1145166666666667, 0.34355, 0.5923587361059505, 0.6855055141740312, 0.7190761834555557, 6.4716856511, 0.1054227272996478, 0.1394517551584083),
    @("Neutro", "L1cam", "Cntn1", "sCs", 3, 1, 6.279227333333334, 18.837682, 0.1779710855497397, 0.2034290786507157, 1, 0.5, 0.0788065, 0.157613, 0.4076412638940496, 0.3144944858259688, 0.4948439288443334, 2.969063573066, 0.07254835825009191, 0.06397732349230739),
    @("sCs", "L1cam", "Cntn1", "FAPs", 2, 1, 8.750795, 17.50159, 0.2480223127622444, 0.1890005537105138, 3, 1, 0.1145166666666667, 0.34355, 0.5923587361059505, 0.6855055141740312, 1.002111874083333, 6.012671244500001, 0.1469181837139178, 0.1295609217505024),
    @("sCs", "L1cam", "Cntn1", "sCs", 2, 1, 8.750795, 17.50159, 0.2480223127622444, 0.1890005537105138, 1, 0.5, 0.0788065, 0.157613, 0.4076412638940496, 0.3144944858259688, 0.6896195261675, 2.75847810467, 0.1011041290483266, 0.05943963196001144)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Count; $c++) {
        $col = $c + 1
        $ws.Cells.Item($row, $col).Value = $vals[$c]
    }
}
